# Atualiza datasets e ajustes das ligas
# Rebuilds the team list in column A/B/C with the new teams inserted in
# their proper (ID-sorted) positions, then re-creates every hyperlink in
# column C so the internal hyperlink relationships stay consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks on the data range first so nothing stale is
# left behind once the cell values below are rewritten/extended.
$ws.Range("C2:C15").Hyperlinks.Delete()

# Final, ID-ordered roster (row, team name, numeric team id).
$rows = @(
    @(2,  "La Primeira Patada Es Nuestra", 32966),
    @(3,  "SC ÉoINTER! ", 184499),
    @(4,  "FBC Colorado", 186283),
    @(5,  "Doug Leal F.C", 287965),
    @(6,  "Texas Club 2026", 1273719),
    @(7,  "C R Juvenal", 1488983),
    @(8,  "JV5 Tricolor Gaúcho", 1747619),
    @(9,  "Medonho´s F.C. ", 1867254),
    @(10, "GaúchoDaFronteira F.C", 2371918),
    @(11, "Esquadrão Gazembrino", 2916559),
    @(12, " NHU PORÃ SAF.", 4088673),
    @(13, "SC 100 Sono", 14709358),
    @(14, "GrioTeam", 14933455),
    @(15, "GE Bebum", 16411206),
    @(16, "bugredasmissões ", 19209079),
    @(17, "Pontaç0 F.C.", 20651178),
    @(18, "lsauer fc", 44810918),
    @(19, "Grêmio_Campeão_LA_27", 47775950)
)

foreach ($r in $rows) {
    $row = $r[0]
    $name = $r[1]
    $id = $r[2]
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $id
    $ws.Range("C$row").Value = "https://cartola.globo.com/#!/time/$id"
}

# Re-create every hyperlink on column C (rows 2..19) in order so the
# relationship ids line up the same way Excel itself would assign them.
foreach ($r in $rows) {
    $row = $r[0]
    $id = $r[2]
    $url = "https://cartola.globo.com/#!/time/$id"
    $ws.Hyperlinks.Add($ws.Range("C$row"), "https://cartola.globo.com/", "!/time/$id", "", $url)
}

# Hyperlinks.Add() stamps a brand-new "applyFont" style onto every cell it
# touches; fold all of column C back onto the same Hyperlink-styled xf that
# the original rows already use (copy style only, keep the text/value).
foreach ($r in $rows) {
    $row = $r[0]
    $id = $r[2]
    $url = "https://cartola.globo.com/#!/time/$id"
    $ws.Range("C2").Copy($ws.Range("C$row"))
    $ws.Range("C$row").Value = $url
}
